$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 update (Processing Time, column J)
$ws.Range("J1").Value = -1512081046

# Row 2 updates (Total Distance / Total Travel Time / Total Cost)
$ws.Range("D2").Value = 9663.8544921875
$ws.Range("F2").Value = 5637.248453776042
$ws.Range("H2").Value = 9663.85

# Row 5 (route) update: new permutation of 1..52 across B5:BA5
$newRoute = @(30,29,47,26,14,52,13,27,28,11,51,12,16,46,44,50,20,23,1,34,35,36,37,48,24,5,25,4,33,43,6,15,38,40,39,49,32,45,19,10,9,8,41,22,31,18,3,21,17,42,7,2)

for ($i = 0; $i -lt $newRoute.Length; $i++) {
    $ws.Cells.Item(5, $i + 2).Value = $newRoute[$i]
}
